$d = $word.ActiveDocument

$d.Content.Find.Execute("branch alternate", $true, $false, $false, $false, $false,
                         $true, 1, $false, "branch main", 2)
